$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("example2")

# Insert a new row above the current row 5 ("b" / "Timeseries b" / "(EUR)" / 10 / 50 / 60),
# pushing it down to row 6, and label it with a new sub-heading in column C.
[void]$ws2.Rows.Item(5).Insert()
$ws2.Range("C5").Value = "Timeseries b in euros"
$ws2.Range("C5").Font.Italic = $true

# Make the "example2" sheet the active tab / selected sheet, with the given selection,
# which also clears "first_sheet"'s tabSelected flag.
[void]$ws2.Activate()
[void]$ws2.Range("D11").Select()
